$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item("Metadata")
$wsElem = $wb.Worksheets.Item("Elements")

# Update metadata Date and Description values
$wsMeta.Range("B8").Value = "2026-01-14T15:34:52+00:00"
$wsMeta.Range("B12").Value = "Instructions au patient"

# The Description string is shared with Elements!M2 (Definition of the root element) -
# update it too so both cells keep referencing the same (renamed) text.
$wsElem.Range("M2").Value = "Instructions au patient"

# Elements sheet: row 5 Short/Definition now reuse the (now-identical) Description text
$wsElem.Range("L5").Value = "Instructions au patient"
$wsElem.Range("M5").Value = "Instructions au patient"
